# Weekly SEBI downloads workbook update ("Circulars and Regulations prompt added")
#
# - New row inserted at row 2: SEBI / Circulars (2025-12-12)
# - Previous row 2 content (AIF / Regulations, 2025-12-11) moved down to row 3
# - New row added at row 4: Listed Companies / Circular-BSE (2025-12-09)
# - Leftover placeholder cells in rows 5, 8, 9, 11, 12 and 13 removed

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain-text cell format used by the rest of the table body (no explicit
# number format, default font). Used later to "stamp" newly written cells
# so they match the existing look instead of picking up a new style.
$plainTpl = $ws.Range("H2")

# ===========================================================================
# Row 2 : SEBI / Circulars  (new top entry)
# ===========================================================================
$ws.Range("C2").NumberFormat = "@"
$ws.Range("E2").NumberFormat = "@"

$ws.Range("A2").Value = "SEBI"
$ws.Range("B2").Value = "Circulars"
$ws.Range("C2").Value = "2025"
$ws.Range("D2").Value = "December"
$ws.Range("E2").Value = "2025-12-12"
$ws.Range("F2").Value = "Provisions relating to Strengthening Governance of Market Infrastructure Institutions (MIIs)"

$g2 = $ws.Range("G2")
$g2.Value = "https://www.sebi.gov.in/sebi_data/attachdocs/dec-2025/1765535283954.pdf"
$ws.Hyperlinks.Add($g2, "https://www.sebi.gov.in/sebi_data/attachdocs/dec-2025/1765535283954.pdf") | Out-Null
# Newly inserted hyperlinks in this sheet render with a plain blue
# underline rather than the workbook's theme hyperlink color.
$g2.Font.Color = 255 * 65536

$ws.Range("H2").Value = "1765535283954.pdf"
$ws.Range("I2").Value = "/Users/admin/Downloads/Tejomaya_pdfs_test/Akshayam Data/SEBI/Circulars/2025/December/1765535283954.pdf"

$plainTpl.Copy()
$ws.Range("A2:F2").PasteSpecial(-4122) | Out-Null
$ws.Range("H2:I2").PasteSpecial(-4122) | Out-Null

# ===========================================================================
# Row 3 : AIF / Regulations (the entry that used to be in row 2)
# ===========================================================================
$ws.Range("C3").NumberFormat = "@"
$ws.Range("E3").NumberFormat = "@"

$ws.Range("A3").Value = "AIF"
$ws.Range("B3").Value = "Regulations"
$ws.Range("C3").Value = "2025"
$ws.Range("D3").Value = "December"
$ws.Range("E3").Value = "2025-12-11"
$ws.Range("F3").Value = "Securities and Exchange Board of India (Real Estate Investment Trusts) Regulations, 2014  [Last amended on December 11, 2025]"

$g3 = $ws.Range("G3")
$g3.Value = "https://www.sebi.gov.in/sebi_data/attachdocs/dec-2025/1765541474113.pdf"
$ws.Hyperlinks.Add($g3, "https://www.sebi.gov.in/sebi_data/attachdocs/dec-2025/1765541474113.pdf") | Out-Null

$ws.Range("H3").Value = "1765541474113.pdf"
$ws.Range("I3").Value = "/Users/admin/Downloads/Tejomaya_pdfs_test/Akshayam Data/AIF/Regulations/2025/December/1765541474113.pdf"

$plainTpl.Copy()
$ws.Range("A3:F3").PasteSpecial(-4122) | Out-Null
$ws.Range("H3:I3").PasteSpecial(-4122) | Out-Null

# ===========================================================================
# Row 4 : Listed Companies / Circular-BSE (brand new entry)
# ===========================================================================
$ws.Range("C4").NumberFormat = "@"
$ws.Range("E4").NumberFormat = "@"

$ws.Range("A4").Value = "Listed Companies"
$ws.Range("B4").Value = "Circular-BSE"
$ws.Range("C4").Value = "2025"
$ws.Range("D4").Value = "December"
$ws.Range("E4").Value = "2025-12-09"
$ws.Range("F4").Value = "Timeline for submission of information by the Issuer to the Debenture Trustee(s)"

$g4 = $ws.Range("G4")
$g4.Value = "https://www.bseindia.com/markets/MarketInfo/DownloadAttach.aspx?id=20251209-53&attachedId=3c58b20b-4d62-43be-b20f-8a824667e02b"
$ws.Hyperlinks.Add($g4, "https://www.bseindia.com/markets/MarketInfo/DownloadAttach.aspx?id=20251209-53&attachedId=3c58b20b-4d62-43be-b20f-8a824667e02b") | Out-Null

$ws.Range("H4").Value = "Timeline for submission of information by the Issuer to DT.pdf"
$ws.Range("I4").Value = "/Users/admin/Downloads/Tejomaya_pdfs_test/Akshayam Data/Listed Companies/Circular-BSE/2025/December/Timeline for submission of information by the Issuer to DT.pdf"

$plainTpl.Copy()
$ws.Range("A4:F4").PasteSpecial(-4122) | Out-Null
$ws.Range("H4:I4").PasteSpecial(-4122) | Out-Null

# ===========================================================================
# Remove the stray placeholder cells that used to pad out the sheet.
# Rows 6 and 7 keep their existing empty, styled G cell untouched.
# ===========================================================================
$ws.Range("G5").Clear() | Out-Null
$ws.Range("G8").Clear() | Out-Null
$ws.Range("G9").Clear() | Out-Null
$ws.Range("G11").Clear() | Out-Null
$ws.Range("G12").Clear() | Out-Null
$ws.Range("G13").Clear() | Out-Null

Write-Output "edit complete"
